$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "50.872.02"
Set-TextValue "E2" "  -16.03%  "
Set-TextValue "D3" "2.234.10"
Set-TextValue "E3" "  -22.98%  "
Set-TextValue "E4" "  +0.30%  "
Set-TextValue "D5" "426.42"
Set-TextValue "E5" "  -18.95%  "
Set-TextValue "D6" "115.84"
Set-TextValue "E6" "  -18.91%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.06%  "
Set-TextValue "D8" "0.451"
Set-TextValue "E8" "  -17.34%  "
Set-TextValue "D9" "2.235.37"
Set-TextValue "E9" "  -23.10%  "
Set-TextValue "D10" "5.12"
Set-TextValue "E10" "  -14.24%  "
Set-TextValue "D11" "0.0839"
Set-TextValue "E11" "  -21.53%  "
Set-TextValue "D12" "0.293"
Set-TextValue "E12" "  -18.11%  "
Set-TextValue "D14" "2.632.49"
Set-TextValue "E14" "  -22.69%  "
Set-TextValue "D15" "51.009.98"
Set-TextValue "E15" "  -15.81%  "
Set-TextValue "D16" "18.27"
Set-TextValue "E16" "  -18.64%  "
Set-TextValue "E17" "  -19.48%  "
Set-TextValue "D18" "2.266.11"
Set-TextValue "E18" "  -22.10%  "
Set-TextValue "D19" "3.84"
Set-TextValue "E19" "  -22.46%  "
Set-TextValue "D20" "289.87"
Set-TextValue "E20" "  -17.19%  "
Set-TextValue "D21" "0.988"
Set-TextValue "E21" "  -1.18%  "
Set-TextValue "D22" "5.67"
Set-TextValue "E22" "  -0.71%  "
Set-TextValue "D23" "8.49"
Set-TextValue "E23" "  -26.55%  "
Set-TextValue "D24" "4.92"
Set-TextValue "E24" "  -24.26%  "
Set-TextValue "D25" "1.00"
Set-TextValue "E25" "  +0.54%  "
Set-TextValue "D26" "52.29"
Set-TextValue "E26" "  -19.21%  "
Set-TextValue "B27" "Polygon"
Set-TextValue "C27" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D27" "0.357"
Set-TextValue "E27" "  -20.81%  "
Set-TextValue "B28" "WrappedeETH"
Set-TextValue "C28" "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D28" "2.367.58"
Set-TextValue "E28" "  -21.83%  "
Set-TextValue "D29" "0.136"
Set-TextValue "E29" "  -23.26%  "
Set-TextValue "D30" "0.999"
Set-TextValue "E30" "  -0.07%  "
Set-TextValue "D31" "6.65"
Set-TextValue "E31" "  -14.68%  "
Set-TextValue "B32" "Monero"
Set-TextValue "C32" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D32" "141.70"
Set-TextValue "E32" "  -6.78%  "
Set-TextValue "B33" "PEPE"
Set-TextValue "C33" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D33" "0.0₃0620"
Set-TextValue "E33" "  -27.32%  "
Set-TextValue "D34" "16.28"
Set-TextValue "E34" "  -16.45%  "
Set-TextValue "D35" "1.28"
Set-TextValue "E35" "  -23.50%  "
Set-TextValue "D36" "4.55"
Set-TextValue "E36" "  -17.95%  "
Set-TextValue "B37" "FirstDigitalUSD"
Set-TextValue "C37" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D37" "0.998"
Set-TextValue "E37" "  +0.03%  "
Set-TextValue "B38" "Fetch.AI"
Set-TextValue "C38" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D38" "0.772"
Set-TextValue "E38" "  -22.07%  "
Set-TextValue "D39" "3.27"
Set-TextValue "E39" "  -23.61%  "
Set-TextValue "D40" "0.960"
Set-TextValue "E40" "  -19.37%  "
Set-TextValue "D41" "31.49"
Set-TextValue "E41" "  -16.20%  "
Set-TextValue "D42" "10.11"
Set-TextValue "E42" "  -2.03%  "
Set-TextValue "D43" "0.543"
Set-TextValue "E43" "  -16.22%  "
Set-TextValue "D44" "0.0486"
Set-TextValue "E44" "  -15.74%  "
Set-TextValue "D45" "3.01"
Set-TextValue "E45" "  -18.67%  "
Set-TextValue "D46" "1.844.75"
Set-TextValue "E46" "  -19.41%  "
Set-TextValue "D47" "1.12"
Set-TextValue "E47" "  -23.29%  "
Set-TextValue "D48" "0.0796"
Set-TextValue "E48" "  -12.94%  "
Set-TextValue "D49" "0.0197"
Set-TextValue "E49" "  -16.62%  "
Set-TextValue "D50" "3.88"
Set-TextValue "E50" "  -21.07%  "
Set-TextValue "B51" "ZEEBU"
Set-TextValue "C51" "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
Set-TextValue "D51" "4.61"
Set-TextValue "E51" "  -5.13%  "
